# agrego atributos a usuario
#
# The document currently has a single (empty) paragraph that carries the
# "_GoBack" bookmark. We need to insert four new "attribute" lines of
# pseudo-code before it:
#   String nombre;
#   int Cedula;
#   date fecha;
#   long cedula;
# The last line ("long cedula;") ends up living in the original paragraph
# (so it keeps the bookmark), while the first three lines become brand new
# paragraphs above it. "String", "int" and "long" are marked as
# spell-check-flagged tokens (w:proofErr spellStart/spellEnd), matching
# what Word's spell checker does for identifiers it doesn't recognize.

$d = $word.ActiveDocument

# The target paragraph is the sole existing paragraph (holds the bookmark).
$target = $d.Paragraphs(1)

# Collapse the range to its very start so the inserted XML lands *before*
# the existing (bookmarked) content instead of replacing it.
$r = $target.Range
$r.Collapse(1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>String</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nombre;</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Cedula;</w:t></w:r></w:p><w:p><w:r><w:t>date fecha;</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>long</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cedula;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
